$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 ("MinCode" item) - relabel the Chinese description and the code name
# to match the QC1297 code-length fields naming convention.
$ws.Range("C15").Value = "代碼最小長度"
$ws.Range("C16").Value = "代碼最大長度"
$ws.Range("B15").Value = "MinCodeLength"
$ws.Range("B16").Value = "MaxCodeLength"

# Leave the cursor where the author last left it before saving.
$ws.Range("G12").Select()
